$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.645.35"
$ws.Range("E2").Value = "  +3.51%  "

$ws.Range("D3").Value = "1.610.27"
$ws.Range("E3").Value = "  +2.85%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.87"
$ws.Range("E5").Value = "  +1.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.519"
$ws.Range("E6").Value = "  +2.17%  "

$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "27.15"
$ws.Range("E8").Value = "  +9.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.64"
$ws.Range("E9").Value = "  -1.26%  "

$ws.Range("E10").Value = "  +2.46%  "

$ws.Range("E11").Value = "  +2.49%  "

$ws.Range("E12").Value = "  +1.41%  "

$ws.Range("D13").Value = "1.839.93"
$ws.Range("E13").Value = "  +2.88%  "

$ws.Range("D14").Value = "1.616.41"
$ws.Range("E14").Value = "  +2.61%  "

$ws.Range("D15").Value = "29.651.81"
$ws.Range("E15").Value = "  +3.46%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.538"
$ws.Range("E16").Value = "  +4.08%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.73"
$ws.Range("E17").Value = "  +2.67%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.58"
$ws.Range("E18").Value = "  +3.44%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.59"
$ws.Range("E19").Value = "  +5.65%  "

$ws.Range("E20").Value = "  +3.83%  "

$ws.Range("D21").Value = "0.0₃0696"
$ws.Range("E21").Value = "  +2.00%  "

$ws.Range("E22").Value = "  +0.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.00"
$ws.Range("E23").Value = "  +1.69%  "

$ws.Range("E24").Value = "  +1.91%  "

$ws.Range("E25").Value = "  +1.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.03"
$ws.Range("E26").Value = "  +2.22%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.33"
$ws.Range("E27").Value = "  +3.74%  "

$ws.Range("E28").Value = "  +2.49%  "

$ws.Range("E29").Value = "  +2.82%  "

$ws.Range("E30").Value = "  +0.14%  "

$ws.Range("E31").Value = "  +3.61%  "

$ws.Range("E32").Value = "  +1.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.23"
$ws.Range("E33").Value = "  +1.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.13"
$ws.Range("E34").Value = "  +4.48%  "

$ws.Range("D35").Value = "1.422.12"
$ws.Range("E35").Value = "  +1.18%  "

$ws.Range("E36").Value = "  -0.39%  "

$ws.Range("E37").Value = "  +4.81%  "

$ws.Range("E38").Value = "  +5.49%  "

$ws.Range("E39").Value = "  +0.44%  "

$ws.Range("E40").Value = "  +2.16%  "

$ws.Range("E41").Value = "  +4.47%  "

$ws.Range("E42").Value = "  +2.13%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "54.45"
$ws.Range("E43").Value = "  +27.41%  "

$ws.Range("E44").Value = "  +6.26%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.801"
$ws.Range("E45").Value = "  +4.22%  "

$ws.Range("E46").Value = "  +0.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "66.06"
$ws.Range("E47").Value = "  +3.30%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.31"
$ws.Range("E48").Value = "  +1.66%  "

$ws.Range("D49").Value = "1.750.55"
$ws.Range("E49").Value = "  +3.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.891"
$ws.Range("E50").Value = "  +3.66%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "86.87"
$ws.Range("E51").Value = "  +2.39%  "
